# Apply the "single work request" report refresh to the Work Report sheet.
# Recomputed totals, refreshed line items (one row removed), and a new
# "Report Generated On" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates -------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"

$ws.Range("C8").Value = 4728.01      # Total Billed Amount
$ws.Range("C9").Value = 9            # Total Line Items
$ws.Range("G10").Value = ""          # Scope ID # cleared

# --- Line item detail rows -----------------------------------------------------
# Row 16: Point 01 / DEG-4-SNA - pricing filled in
$ws.Range("H16").Value = 572.5

# Row 17: Point 04 / PLA-DLOC - pricing filled in
$ws.Range("H17").Value = 952.8

# Row 18: now CON-2-AAI-3-P (was SVW-2-TP-CUS-CC)
$ws.Range("B18").Value = "CON-2-AAI-3-P"
$ws.Range("D18").Value = "CON,#2 AWG,AA Corr,Three,Poly"
$ws.Range("F18").Value = 30
$ws.Range("H18").Value = 21.6

# Row 19: now CNC-NTI-10 (was CON-2-AAI-3-P)
$ws.Range("B19").Value = "CNC-NTI-10"
$ws.Range("D19").Value = "CNC,splice Non-Tension Insul,336-1033"
$ws.Range("E19").Value = "EA"
$ws.Range("F19").Value = 6
$ws.Range("H19").Value = 121.74

# Row 20: now Point 09 / PLA-DLOC (was Point 01 / CNC-NTI-10)
$ws.Range("A20").Value = "Point 09"
$ws.Range("B20").Value = "PLA-DLOC"
$ws.Range("D20").Value = "PLA,Difficult Location"
$ws.Range("H20").Value = 714.6

# Row 21: now Point 04 (was Point 09), PLA-DLOC unchanged
$ws.Range("A21").Value = "Point 04"
$ws.Range("H21").Value = 714.6

# Row 22: now Point 05 (was Point 04), PLA-DLOC unchanged
$ws.Range("A22").Value = "Point 05"
$ws.Range("H22").Value = 714.6

# Row 23: now Point 06 / CNC-HTA-40 (was Point 05 / PLA-DLOC)
$ws.Range("A23").Value = "Point 06"
$ws.Range("B23").Value = "CNC-HTA-40"
$ws.Range("D23").Value = "Compression Connector H-Tap Assembly 4/0"
$ws.Range("F23").Value = 21
$ws.Range("H23").Value = 200.97

# Row 24: now Point 08 / PLA-DLOC (was Point 06 / CNC-HTA-40)
$ws.Range("A24").Value = "Point 08"
$ws.Range("B24").Value = "PLA-DLOC"
$ws.Range("D24").Value = "PLA,Difficult Location"
$ws.Range("F24").Value = 6
$ws.Range("H24").Value = 714.6

# Row 25 (old "Point 08 / PLA-DLOC" line item) is dropped entirely - deleting
# it shifts the old TOTAL row (26) up to row 25, carrying its styles/merge
# along with it.
$ws.Rows.Item(25).Delete()

# --- Recalculated TOTAL row (now row 25) ---------------------------------------
$ws.Range("H25").Value = 4728.009999999999
